$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 763.7143
$ws.Range("I8").Value = 57.666668
$ws.Range("J8").Value = 5000
$ws.Range("K8").Value = 173.000004
$ws.Range("L8").Value = 15000
$ws.Range("M8").Value = -34.00000399999999
$ws.Range("N8").Value = -15278
$ws.Range("H33").Value = 548.0952
$ws.Range("I33").Value = 115.14286
$ws.Range("K33").Value = 115.14286
$ws.Range("M33").Value = 113.85714
$ws.Range("H70").Value = 4814.7144
$ws.Range("I70").Value = 500
$ws.Range("J70").Value = 6540.6
$ws.Range("K70").Value = 1500
$ws.Range("L70").Value = 19621.8
$ws.Range("M70").Value = -1230
$ws.Range("N70").Value = -20161.8
$ws.Range("H73").Value = 4814.7144
$ws.Range("I73").Value = 500
$ws.Range("J73").Value = 6540.6
$ws.Range("K73").Value = 1500
$ws.Range("L73").Value = 19621.8
$ws.Range("M73").Value = -564
$ws.Range("N73").Value = -21493.8
$ws.Range("H132").Value = 654797.2
$ws.Range("I132").Value = 1282.1846
$ws.Range("J132").Value = 4902645
$ws.Range("K132").Value = 3846.5538
$ws.Range("L132").Value = 14707935
$ws.Range("M132").Value = -1316.5538
$ws.Range("N132").Value = -14712995
$ws.Range("H137").Value = 1493742.9
$ws.Range("I137").Value = 2084271.2
$ws.Range("J137").Value = 1881.579
$ws.Range("K137").Value = 6252813.6
$ws.Range("L137").Value = 5644.737
$ws.Range("M137").Value = -6250263.6
$ws.Range("N137").Value = -10744.737
$ws.Range("H138").Value = 1833230.1
$ws.Range("I138").Value = 874.7692
$ws.Range("J138").Value = 3207496.8
$ws.Range("K138").Value = 2624.3076
$ws.Range("L138").Value = 9622490.399999999
$ws.Range("M138").Value = 2515.6924
$ws.Range("N138").Value = -9632770.399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1297.54
$ws.Range("I32").Value = 852.4125
$ws.Range("J32").Value = 3078.05
$ws.Range("K32").Value = 852.4125
$ws.Range("L32").Value = 3078.05
$ws.Range("M32").Value = -565.4125
$ws.Range("N32").Value = -3652.05
$ws.Range("H45").Value = 1560
$ws.Range("I45").Value = 1163.6364
$ws.Range("K45").Value = 1163.6364
$ws.Range("M45").Value = -786.6364000000001
$ws.Range("H88").Value = 4768.4165
$ws.Range("I88").Value = 2214.2856
$ws.Range("J88").Value = 5820.1177
$ws.Range("K88").Value = 2214.2856
$ws.Range("L88").Value = 5820.1177
$ws.Range("M88").Value = -1808.2856
$ws.Range("N88").Value = -6632.1177
$ws.Range("H91").Value = 4768.4165
$ws.Range("I91").Value = 2214.2856
$ws.Range("J91").Value = 5820.1177
$ws.Range("K91").Value = 2214.2856
$ws.Range("L91").Value = 5820.1177
$ws.Range("M91").Value = -810.2856000000002
$ws.Range("N91").Value = -8628.117699999999
$ws.Range("H122").Value = 5292884.5
$ws.Range("I122").Value = 2148.7144
$ws.Range("J122").Value = 15874356
$ws.Range("K122").Value = 6446.1432
$ws.Range("L122").Value = 47623068
$ws.Range("M122").Value = -3996.1432
$ws.Range("N122").Value = -47627968
$ws.Range("H139").Value = 47400
$ws.Range("J139").Value = 47400
$ws.Range("L139").Value = 47400
$ws.Range("N139").Value = -57680

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 10243.286
$ws.Range("I86").Value = 11808.608
$ws.Range("J86").Value = 3042.8
$ws.Range("K86").Value = 11808.608
$ws.Range("L86").Value = 3042.8
$ws.Range("M86").Value = -10685.608
$ws.Range("N86").Value = -5288.8
$ws.Range("H89").Value = 10243.286
$ws.Range("I89").Value = 11808.608
$ws.Range("J89").Value = 3042.8
$ws.Range("K89").Value = 59043.04
$ws.Range("L89").Value = 15214
$ws.Range("M89").Value = -53427.04
$ws.Range("N89").Value = -26446
$ws.Range("H141").Value = 45251
$ws.Range("J141").Value = 46765
$ws.Range("L141").Value = 46765
$ws.Range("N141").Value = -57125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 1799.8
$ws.Range("J14").Value = 2000
$ws.Range("L14").Value = 2000
$ws.Range("N14").Value = -2340
$ws.Range("H99").Value = 10500
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 10500
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 10500
$ws.Range("N99").Value = -13496
$ws.Range("M99").ClearContents()
$ws.Range("H126").Value = 10500
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 10500
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 31500
$ws.Range("N126").Value = -36440
$ws.Range("M126").ClearContents()
$ws.Range("H134").Value = 20342.16
$ws.Range("I134").Value = 967.4
$ws.Range("J134").Value = 99602.55
$ws.Range("K134").Value = 2902.2
$ws.Range("L134").Value = 298807.65
$ws.Range("M134").Value = -367.1999999999998
$ws.Range("N134").Value = -303877.65

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 48899.348
$ws.Range("I70").Value = 88198.75
$ws.Range("J70").Value = 6027.273
$ws.Range("K70").Value = 88198.75
$ws.Range("L70").Value = 6027.273
$ws.Range("M70").Value = -87928.75
$ws.Range("N70").Value = -6567.273
$ws.Range("H73").Value = 48899.348
$ws.Range("I73").Value = 88198.75
$ws.Range("J73").Value = 6027.273
$ws.Range("K73").Value = 88198.75
$ws.Range("L73").Value = 6027.273
$ws.Range("M73").Value = -87262.75
$ws.Range("N73").Value = -7899.273
$ws.Range("H102").Value = 1080
$ws.Range("I102").Value = 1000
$ws.Range("J102").Value = 1133.3334
$ws.Range("K102").Value = 1000
$ws.Range("L102").Value = 1133.3334
$ws.Range("M102").Value = 622
$ws.Range("N102").Value = -4377.3334
$ws.Range("H126").Value = 1675
$ws.Range("I126").Value = 1410
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 4230
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -1760
$ws.Range("N126").Value = -13940
$ws.Range("H127").Value = 49800
$ws.Range("J127").Value = 49800
$ws.Range("L127").Value = 49800
$ws.Range("N127").Value = -59720
$ws.Range("H132").Value = 48959.383
$ws.Range("I132").Value = 29084.111
$ws.Range("J132").Value = 168211
$ws.Range("K132").Value = 87252.333
$ws.Range("L132").Value = 504633
$ws.Range("M132").Value = -84722.333
$ws.Range("N132").Value = -509693
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H138").Value = 53214.5
$ws.Range("J138").Value = 53214.5
$ws.Range("L138").Value = 53214.5
$ws.Range("N138").Value = -63494.5
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("M140").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2338.3076
$ws.Range("I7").Value = 2233.3333
$ws.Range("K7").Value = 2233.3333
$ws.Range("M7").Value = -2121.3333
$ws.Range("H16").Value = 4475.95
$ws.Range("I16").Value = 1108.5
$ws.Range("J16").Value = 12333.333
$ws.Range("K16").Value = 1108.5
$ws.Range("L16").Value = 12333.333
$ws.Range("M16").Value = -938.5
$ws.Range("N16").Value = -12673.333
$ws.Range("H40").Value = 2986.7856
$ws.Range("I40").Value = 2665
$ws.Range("J40").Value = 4166.6665
$ws.Range("K40").Value = 2665
$ws.Range("L40").Value = 4166.6665
$ws.Range("M40").Value = -2529
$ws.Range("N40").Value = -4438.6665
$ws.Range("H61").Value = 2917.15
$ws.Range("I61").Value = 2874.5
$ws.Range("J61").Value = 3016.6667
$ws.Range("K61").Value = 2874.5
$ws.Range("L61").Value = 3016.6667
$ws.Range("M61").Value = -2672.5
$ws.Range("N61").Value = -3420.6667
$ws.Range("H93").Value = 1208.6364
$ws.Range("I93").Value = 1149.5
$ws.Range("J93").Value = 1366.3334
$ws.Range("K93").Value = 1149.5
$ws.Range("L93").Value = 1366.3334
$ws.Range("M93").Value = 98.5
$ws.Range("N93").Value = -3862.3334
$ws.Range("H113").Value = 2917.15
$ws.Range("I113").Value = 2874.5
$ws.Range("J113").Value = 3016.6667
$ws.Range("K113").Value = 2874.5
$ws.Range("L113").Value = 3016.6667
$ws.Range("M113").Value = -704.5
$ws.Range("N113").Value = -7356.6667
$ws.Range("H122").Value = 3053.037
$ws.Range("I122").Value = 2261.9285
$ws.Range("J122").Value = 3905
$ws.Range("K122").Value = 6785.7855
$ws.Range("L122").Value = 11715
$ws.Range("M122").Value = -4335.7855
$ws.Range("N122").Value = -16615
$ws.Range("H126").Value = 2338.3076
$ws.Range("I126").Value = 2233.3333
$ws.Range("K126").Value = 6699.999899999999
$ws.Range("M126").Value = -4229.999899999999
$ws.Range("H132").Value = 35587.785
$ws.Range("I132").Value = 25736.162
$ws.Range("J132").Value = 59122.223
$ws.Range("K132").Value = 77208.486
$ws.Range("L132").Value = 177366.669
$ws.Range("M132").Value = -74678.486
$ws.Range("N132").Value = -182426.669
$ws.Range("H136").Value = 70229.25999999999
$ws.Range("I136").Value = 55778.4
$ws.Range("J136").Value = 96503.55
$ws.Range("K136").Value = 167335.2
$ws.Range("L136").Value = 289510.65
$ws.Range("M136").Value = -164785.2
$ws.Range("N136").Value = -294610.65
$ws.Range("H137").Value = 29960.6
$ws.Range("J137").Value = 31433
$ws.Range("L137").Value = 31433
$ws.Range("N137").Value = -41633
$ws.Range("H139").Value = 43987.5
$ws.Range("J139").Value = 43987.5
$ws.Range("L139").Value = 43987.5
$ws.Range("N139").Value = -54267.5
$ws.Range("H141").Value = 49715
$ws.Range("J141").Value = 49715
$ws.Range("L141").Value = 49715
$ws.Range("N141").Value = -60075

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 54000
$ws.Range("J95").Value = 54000
$ws.Range("L95").Value = 54000
$ws.Range("N95").Value = -59492
$ws.Range("H122").Value = 3618.818
$ws.Range("I122").Value = 2852
$ws.Range("J122").Value = 4057
$ws.Range("K122").Value = 8556
$ws.Range("L122").Value = 12171
$ws.Range("M122").Value = -6106
$ws.Range("N122").Value = -17071
$ws.Range("H126").Value = 1169.6428
$ws.Range("I126").Value = 1175.3846
$ws.Range("J126").Value = 1095
$ws.Range("K126").Value = 3526.1538
$ws.Range("L126").Value = 3285
$ws.Range("M126").Value = -1056.1538
$ws.Range("N126").Value = -8225
$ws.Range("H132").Value = 41415.49
$ws.Range("I132").Value = 33407.246
$ws.Range("J132").Value = 73982.336
$ws.Range("K132").Value = 100221.738
$ws.Range("L132").Value = 221947.008
$ws.Range("M132").Value = -97691.738
$ws.Range("N132").Value = -227007.008
$ws.Range("H136").Value = 45680.844
$ws.Range("I136").Value = 32873.484
$ws.Range("J136").Value = 74040
$ws.Range("K136").Value = 98620.452
$ws.Range("L136").Value = 222120
$ws.Range("M136").Value = -96070.452
$ws.Range("N136").Value = -227220
